# ST3 Tool Monitoring sheet update
# - Header row (row 1) is extended from column W out to column AU with new
#   "Cap Handle Oiling" / "Cap Handle and Threaded ferrule Assembly" checks.
# - Remark / QA-Sign / Engg-Sign headers shift from U/V/W to AS/AT/AU.
# - Data row (row 2) gets new sample values, including 24 brand-new
#   OBS/cnt columns, and the trailing Remark/QA-Sign/Engg-Sign values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- New header text for every column A1:AU1 ----
$headers = [ordered]@{
    "A"  = "Date"
    "B"  = "Shift"
    "C"  = "Cap Handle Oiling_TF_Check Fixture_OBS"
    "D"  = "Cap Handle Oiling_TF_Check Fixture_Tool_Life_cnt"
    "E"  = "Cap Handle Oiling_TF_Check the Tip of the Fixture_OBS"
    "F"  = "Cap Handle Oiling_TF_Check the Tip of the Fixture_Life_cnt"
    "G"  = "Cap Handle Oiling_TF_Check Fixture wear out_OBS"
    "H"  = "Cap Handle Oiling_TF_Check Fixture wear out_Tool_Life_cnt"
    "I"  = "Cap Handle Oiling_BF_Check Fixture_OBS"
    "J"  = "Cap Handle Oiling_BF_Check Fixture_Tool_Life_cnt"
    "K"  = "Cap Handle Oiling_BF_Check the Tip of the Fixture_OBS"
    "L"  = "Cap Handle Oiling_BF_Check the Tip of the Fixture_Tool_Life_cnt"
    "M"  = "Cap Handle Oiling_BF_Check Fixture wear out_OBS"
    "N"  = "Cap Handle Oiling_BF_Check Fixture wear out_Tool_Life_cnt"
    "O"  = "Cap Handle Oiling_SN_Check Nozzle_OBS"
    "P"  = "Cap Handle Oiling_SN_Check Nozzle_Tool_Life_cnt"
    "Q"  = "Cap Handle Oiling_SN_Check the Holes of the Nozzle_OBS"
    "R"  = "Cap Handle Oiling_SN_Check the Holes of the Nozzle_Tool_Life_cnt"
    "S"  = "Cap Handle Oiling_SN_Check Nozzle wear out_OBS"
    "T"  = "Cap Handle Oiling_SN_Check Nozzle wear out_Tool_Life_cnt"
    "U"  = "Cap Handle and Threaded ferrule Assembly_BF_Check Fixture_OBS"
    "V"  = "Cap Handle and Threaded ferrule Assembly_BF_Check Fixture_Tool_Life_cnt"
    "W"  = "Cap Handle and Threaded ferrule Assembly_BF_Check the Tip of the Fixture_OBS"
    "X"  = "Cap Handle and Threaded ferrule Assembly_BF_Check the Tip of the Fixture_Tool_Life_cnt"
    "Y"  = "Cap Handle and Threaded ferrule Assembly_BF_Check Fixture wear out_OBS"
    "Z"  = "Cap Handle and Threaded ferrule Assembly_BF_Check Fixture wear out_Tool_Life_cnt"
    "AA" = "Cap Handle and Threaded ferrule Assembly_SPD_Check the Dolly_OBS"
    "AB" = "Cap Handle and Threaded ferrule Assembly_SPD_Check the Dolly_Tool_Life_cnt"
    "AC" = "Cap Handle and Threaded ferrule Assembly_SPD_Check the Tip of the Dolly_OBS"
    "AD" = "Cap Handle and Threaded ferrule Assembly_SPD_Check the Tip of the Dolly_cnt"
    "AE" = "Cap Handle and Threaded ferrule Assembly_SPD_Check Dolly wear out_OBS"
    "AF" = "Cap Handle and Threaded ferrule Assembly_SPD_Check Dolly wear out_Tool_Life_cnt"
    "AG" = "Cap Handle and Threaded ferrule Assembly_PF_Check Fixture_OBS"
    "AH" = "Cap Handle and Threaded ferrule Assembly_PF_Check Fixture_Tool_Life_cnt"
    "AI" = "Cap Handle and Threaded ferrule Assembly_PF_Check Fixture of the Fixture_OBS"
    "AJ" = "Cap Handle and Threaded ferrule Assembly_PF_Check Fixture of the Fixture_Tool_Life_cnt"
    "AK" = "Cap Handle and Threaded ferrule Assembly_PF_Check Fixture wear out_OBS"
    "AL" = "Cap Handle and Threaded ferrule Assembly_PF_Check Fixture wear out_Tool_Life_cnt"
    "AM" = "Cap Handle and Threaded ferrule Assembly_FCG_Check Gripper_OBS"
    "AN" = "Cap Handle and Threaded ferrule Assembly_FCG_Check Gripper_Tool_Life_cnt"
    "AO" = "Cap Handle and Threaded ferrule Assembly_FCG_Check the Tip of the Gripper_OBS"
    "AP" = "Cap Handle and Threaded ferrule Assembly_FCG_Check the Tip of the Gripper_Tool_Life_cnt"
    "AQ" = "Cap Handle and Threaded ferrule Assembly_FCG_Check Gripper wear out_OBS"
    "AR" = "Cap Handle and Threaded ferrule Assembly_FCG_Check Gripper wear out_Tool_Life_cnt"
    "AS" = "Remark"
    "AT" = "QA-Sign"
    "AU" = "Engg-Sign"
}

# ---- New data values for every column A2:AU2 ----
$data = [ordered]@{
    "A"  = "2025-02-06T12:57"
    "B"  = "SHIFT3"
    "C"  = "OK"
    "D"  = "1"
    "E"  = "OK"
    "F"  = "2"
    "G"  = "OK"
    "H"  = "3"
    "I"  = "OK"
    "J"  = "4"
    "K"  = "OK"
    "L"  = "5"
    "M"  = "OK"
    "N"  = "6"
    "O"  = "OK"
    "P"  = "7"
    "Q"  = "OK"
    "R"  = "8"
    "S"  = "OK"
    "T"  = "9"
    "U"  = "OK"
    "V"  = "10"
    "W"  = "OK"
    "X"  = "11"
    "Y"  = "OK"
    "Z"  = "12"
    "AA" = "OK"
    "AB" = "13"
    "AC" = "OK"
    "AD" = "14"
    "AE" = "OK"
    "AF" = "15"
    "AG" = "OK"
    "AH" = "16"
    "AI" = "OK"
    "AJ" = "17"
    "AK" = "OK"
    "AL" = "18"
    "AM" = "OK"
    "AN" = "19"
    "AO" = "OK"
    "AP" = "20"
    "AQ" = "OK"
    "AR" = "21"
    "AS" = "000"
    "AT" = "111"
    "AU" = "222"
}

# The style used by the bold/bordered header band lives on A1 already;
# use it as the template for every header cell (A1:T1 already have it, and
# copying it in is a harmless no-op there; U1:AU1 are brand new cells that
# need it applied explicitly).
$headerStyleSource = $ws.Range("A1")

foreach ($col in $headers.Keys) {
    $cell = $ws.Range($col + "1")
    $headerStyleSource.Copy($cell)
    $cell.Value = $headers[$col]
}

foreach ($col in $data.Keys) {
    $cell = $ws.Range($col + "2")
    # Force text storage (so values like "000", "1", "111" keep their exact
    # literal text instead of being auto-coerced into numbers), then strip
    # the temporary text format back off so the cell stays style-less,
    # matching the original unstyled data row.
    $cell.NumberFormat = "@"
    $cell.Value = $data[$col]
    $cell.ClearFormats()
}

Write-Output "ST3 Tool Monitoring sheet updated: header extended to AU1, data row refreshed through AU2."
